$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

function Set-ParagraphText($textRange, $paraIndex, $newText) {
    $para = $textRange.Paragraphs($paraIndex, 1)
    $len = $para.Text.Length
    $sub = $para.Characters(1, $len)
    $sub.Text = $newText
}

# --- Update the three "Tutoring" course paragraphs (a.a. XX/YY -> YYYY/YYYY) ---
Set-ParagraphText $tr 2 "95631 - MACHINE LEARNING AND DATA MINING - 6 cfu (2023/2024)"
Set-ParagraphText $tr 3 "95631 - MACHINE LEARNING AND DATA MINING - 6 cfu (2024/2025)"
Set-ParagraphText $tr 4 "95631 - MACHINE LEARNING AND DATA MINING - 6 cfu (2025/2026)"

# --- Add a new paragraph after "EDBT/ICDT 2024 Joint Conference" ---
$tr.InsertAfter("`rSEBD 2024, 2025")
